$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1111, shifting existing rows 1111-1230 down to 1112-1231
$ws.Rows("1111:1111").Insert()

# Populate the newly inserted row with the new daily price entry
$ws.Range("A1111").Value = 3
$ws.Range("B1111").Value = "Femacal de La Calera"
$ws.Range("C1111").Value = "Coquimbo"
$ws.Range("D1111").Value = 45194
$ws.Range("E1111").Value = 5
$ws.Range("F1111").Value = 100112023
$ws.Range("G1111").Value = "Brócoli"
$ws.Range("H1111").Value = "Sin especificar"
$ws.Range("I1111").Value = "Primera"
$ws.Range("J1111").Value = 2900
$ws.Range("K1111").Value = 700
$ws.Range("L1111").Value = 750
$ws.Range("M1111").Value = 731
$ws.Range("N1111").Value = "`$/unidad"
$ws.Range("O1111").Value = "Provincia de Quillota"
$ws.Range("P1111").Value = 731
$ws.Range("Q1111").Value = 1
$ws.Range("R1111").Value = "Hortaliza"
